$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 36 (shifts existing rows 36-39 down to 37-40)
$ws.Rows.Item(36).Insert()

# Populate the newly inserted row 36 with the new weekly entry
$ws.Cells.Item(36, 1).Value = 5
$ws.Cells.Item(36, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(36, 3).Value = "Maule"
$ws.Cells.Item(36, 4).Value = 44753
$ws.Cells.Item(36, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(36, 5).Value = 7
$ws.Cells.Item(36, 6).Value = 100112043
$ws.Cells.Item(36, 7).Value = "Pepino dulce"
$ws.Cells.Item(36, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(36, 9).Value = "Primera"
$ws.Cells.Item(36, 10).Value = 300
$ws.Cells.Item(36, 11).Value = 15000
$ws.Cells.Item(36, 12).Value = 15000
$ws.Cells.Item(36, 13).Value = 15000
$ws.Cells.Item(36, 14).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(36, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(36, 16).Value = 833
$ws.Cells.Item(36, 17).Value = 18
$ws.Cells.Item(36, 18).Value = "Hortaliza"
